# Auto-generated Excel COM-interop script to apply scheduled-runner market data updates
# across the Moogle_Profits workbook (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets).

$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 33
$ws.Range("H33").Value = 3054.5454
$ws.Range("I33").Value = 700.25
$ws.Range("K33").Value = 700.25
$ws.Range("M33").Value = -471.25
# Row 43
$ws.Range("H43").Value = 2868
$ws.Range("I43").Value = 2521.65
$ws.Range("J43").Value = 4599.75
$ws.Range("K43").Value = 2521.65
$ws.Range("L43").Value = 4599.75
$ws.Range("M43").Value = -2452.65
$ws.Range("N43").Value = -4737.75
# Row 107
$ws.Range("H107").Value = 418.5
$ws.Range("I107").Value = 391.96155
$ws.Range("J107").Value = 533.5
$ws.Range("K107").Value = 391.96155
$ws.Range("L107").Value = 533.5
$ws.Range("M107").Value = 1528.03845
$ws.Range("N107").Value = -4373.5
# Row 113
$ws.Range("H113").Value = 1474817.5
$ws.Range("I113").Value = 5003052.5
$ws.Range("J113").Value = 4719.5835
$ws.Range("K113").Value = 5003052.5
$ws.Range("L113").Value = 4719.5835
$ws.Range("M113").Value = -4999798.5
$ws.Range("N113").Value = -11227.5835
# Row 121
$ws.Range("H121").Value = 3351.2856
$ws.Range("J121").Value = 3351.2856
$ws.Range("L121").Value = 10053.8568
$ws.Range("N121").Value = -13547.8568
# Row 132
$ws.Range("H132").Value = 2545.5833
$ws.Range("I132").Value = 2607.1177
$ws.Range("J132").Value = 1499.5
$ws.Range("K132").Value = 7821.353099999999
$ws.Range("L132").Value = 4498.5
$ws.Range("M132").Value = -5291.353099999999
$ws.Range("N132").Value = -9558.5
# Row 137
$ws.Range("H137").Value = 2452.5715
$ws.Range("I137").Value = 2313.4546
$ws.Range("J137").Value = 2962.6667
$ws.Range("K137").Value = 6940.3638
$ws.Range("L137").Value = 8888.000100000001
$ws.Range("M137").Value = -4390.3638
$ws.Range("N137").Value = -13988.0001
# Row 140
$ws.Range("H140").Value = 69991.664
$ws.Range("J140").Value = 69991.664
$ws.Range("L140").Value = 69991.664
$ws.Range("N140").Value = -80351.664

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 12309.91
$ws.Range("I32").Value = 8150.9165
$ws.Range("K32").Value = 8150.9165
$ws.Range("M32").Value = -7863.9165
# Row 63
$ws.Range("H63").Value = 1262.9
$ws.Range("I63").Value = 1003.2222
$ws.Range("K63").Value = 1003.2222
$ws.Range("M63").Value = -317.2222
# Row 66
$ws.Range("H66").Value = 1262.9
$ws.Range("I66").Value = 1003.2222
$ws.Range("K66").Value = 5016.111
$ws.Range("M66").Value = -1584.111

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 2
$ws.Range("H2").Value = 50000
$ws.Range("J2").Value = 50000
$ws.Range("L2").Value = 50000
$ws.Range("N2").Value = -50226
# Row 22
$ws.Range("H22").Value = 885.9375
$ws.Range("I22").Value = 793.5
$ws.Range("J22").Value = 1040
$ws.Range("K22").Value = 793.5
$ws.Range("L22").Value = 1040
$ws.Range("M22").Value = -620.5
$ws.Range("N22").Value = -1386
# Row 94
$ws.Range("H94").Value = 1520.4286
$ws.Range("I94").Value = 755.0833
$ws.Range("J94").Value = 2540.889
$ws.Range("K94").Value = 755.0833
$ws.Range("L94").Value = 2540.889
$ws.Range("M94").Value = -304.0833
$ws.Range("N94").Value = -3442.889
# Row 134
$ws.Range("H134").Value = 6428.6943
$ws.Range("I134").Value = 4416
$ws.Range("J134").Value = 11003
$ws.Range("K134").Value = 13248
$ws.Range("L134").Value = 33009
$ws.Range("M134").Value = -10713
$ws.Range("N134").Value = -38079

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 13173.571
$ws.Range("I31").Value = 5130.8
$ws.Range("J31").Value = 20485.182
$ws.Range("K31").Value = 5130.8
$ws.Range("L31").Value = 20485.182
$ws.Range("M31").Value = -4835.8
$ws.Range("N31").Value = -21075.182
# Row 34
$ws.Range("H34").Value = 13173.571
$ws.Range("I34").Value = 5130.8
$ws.Range("J34").Value = 20485.182
$ws.Range("K34").Value = 5130.8
$ws.Range("L34").Value = 20485.182
$ws.Range("M34").Value = -4928.8
$ws.Range("N34").Value = -20889.182
# Row 51
$ws.Range("H51").Value = 37409.332
$ws.Range("I51").Value = 24228
$ws.Range("J51").Value = 44000
$ws.Range("K51").Value = 24228
$ws.Range("L51").Value = 44000
$ws.Range("M51").Value = -23492
$ws.Range("N51").Value = -45472
# Row 52
$ws.Range("H52").Value = 60999
$ws.Range("J52").Value = 60999
$ws.Range("L52").Value = 60999
$ws.Range("N52").Value = -61587
# Row 59
$ws.Range("H59").Value = 81166.336
$ws.Range("I59").Value = 35000
$ws.Range("J59").Value = 104249.5
$ws.Range("K59").Value = 35000
$ws.Range("L59").Value = 104249.5
$ws.Range("M59").Value = -33855
$ws.Range("N59").Value = -106539.5
# Row 61
$ws.Range("H61").Value = 37409.332
$ws.Range("I61").Value = 24228
$ws.Range("J61").Value = 44000
$ws.Range("K61").Value = 24228
$ws.Range("L61").Value = 44000
$ws.Range("M61").Value = -23880
$ws.Range("N61").Value = -44696
# Row 95
$ws.Range("H95").Value = 39639.285
$ws.Range("J95").Value = 39639.285
$ws.Range("L95").Value = 39639.285
$ws.Range("N95").Value = -45131.285
# Row 122
$ws.Range("H122").Value = 3090.0715
$ws.Range("I122").Value = 2306.389
$ws.Range("J122").Value = 4500.7
$ws.Range("K122").Value = 6919.167
$ws.Range("L122").Value = 13502.1
$ws.Range("M122").Value = -4469.167
$ws.Range("N122").Value = -18402.1
# Row 132
$ws.Range("H132").Value = 3244.1226
$ws.Range("I132").Value = 3076.4092
$ws.Range("K132").Value = 9229.2276
$ws.Range("M132").Value = -6699.2276
# Row 140
$ws.Range("H140").Value = 87955.39999999999
$ws.Range("J140").Value = 87955.39999999999
$ws.Range("L140").Value = 87955.39999999999
$ws.Range("N140").Value = -98315.39999999999

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 68
$ws.Range("H68").Value = 4920.579
$ws.Range("I68").Value = 5887.7334
$ws.Range("K68").Value = 17663.2002
$ws.Range("M68").Value = -16852.2002
# Row 71
$ws.Range("H71").Value = 4920.579
$ws.Range("I71").Value = 5887.7334
$ws.Range("K71").Value = 52989.6006
$ws.Range("M71").Value = -48933.6006
# Row 107
$ws.Range("H107").Value = 465.28
$ws.Range("J107").Value = 557.125
$ws.Range("L107").Value = 1671.375
$ws.Range("N107").Value = -5511.375
# Row 113
$ws.Range("H113").Value = 1580.1666
$ws.Range("J113").Value = 1695.5
$ws.Range("L113").Value = 5086.5
$ws.Range("N113").Value = -9426.5
# Row 121
$ws.Range("H121").Value = 4548242
$ws.Range("I121").Value = 1366.6666
$ws.Range("J121").Value = 5788299
$ws.Range("K121").Value = 4099.9998
$ws.Range("L121").Value = 17364897
$ws.Range("M121").Value = -2789.9998
$ws.Range("N121").Value = -17367517

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 100
$ws.Range("H100").Value = 94877.5
$ws.Range("J100").Value = 94877.5
$ws.Range("L100").Value = 94877.5
$ws.Range("N100").Value = -97041.5
# Row 132
$ws.Range("H132").Value = 7323.676
$ws.Range("J132").Value = 8736.799999999999
$ws.Range("L132").Value = 26210.4
$ws.Range("N132").Value = -31270.4
# Row 136
$ws.Range("H136").Value = 28639.5
$ws.Range("J136").Value = 28639.5
$ws.Range("L136").Value = 85918.5
$ws.Range("N136").Value = -91018.5

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 55
$ws.Range("H55").Value = 461.13333
$ws.Range("J55").Value = 702.6923
$ws.Range("L55").Value = 702.6923
$ws.Range("N55").Value = -1048.6923
# Row 114
$ws.Range("H114").Value = 99500
$ws.Range("J114").Value = 99500
$ws.Range("L114").Value = 99500
$ws.Range("N114").Value = -108178
# Row 132
$ws.Range("H132").Value = 4330.1113
$ws.Range("J132").Value = 5963.6
$ws.Range("L132").Value = 17890.8
$ws.Range("N132").Value = -22950.8

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 96
$ws.Range("H96").Value = 6607.75
$ws.Range("I96").Value = 2323.8462
$ws.Range("K96").Value = 2323.8462
$ws.Range("M96").Value = -950.8462
# Row 120
$ws.Range("H120").Value = 47125
$ws.Range("J120").Value = 47125
$ws.Range("L120").Value = 47125
$ws.Range("N120").Value = -56801
# Row 132
$ws.Range("H132").Value = 2683.754
$ws.Range("I132").Value = 2548.1133
$ws.Range("J132").Value = 3282.8333
$ws.Range("K132").Value = 7644.3399
$ws.Range("L132").Value = 9848.499899999999
$ws.Range("M132").Value = -5114.3399
$ws.Range("N132").Value = -14908.4999
# Row 135
$ws.Range("H135").Value = 72666.336
$ws.Range("J135").Value = 72666.336
$ws.Range("L135").Value = 72666.336
$ws.Range("N135").Value = -82806.336
# Row 136
$ws.Range("H136").Value = 5857.9395
$ws.Range("I136").Value = 4897
$ws.Range("J136").Value = 10182.167
$ws.Range("K136").Value = 14691
$ws.Range("L136").Value = 30546.501
$ws.Range("M136").Value = -12141
$ws.Range("N136").Value = -35646.501
